$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure the Price/Volume columns keep their original text formatting
# (values are stored as literal text, not numbers/percentages) by forcing
# the cell number format to Text before assigning the new value.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "307.88"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "-4.67%"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "49.11"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "-1.12%"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "5.174"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "-3.26%"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.07731"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "-5.26%"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "4.508"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "-2.04%"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.339"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "13.69%"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "-7.10%"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.1228"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "-9.01%"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.1936"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "-1.18%"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.04673"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "3.01%"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.09330"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "-2.89%"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.1047"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "-0.08%"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.001263"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "-4.93%"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.04179"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "-2.89%"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.005811"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "-2.33%"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.328"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "-2.09%"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.275"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "-6.54%"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "2.77%"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.992"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "-2.35%"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "-5.58%"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.3039"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "-0.41%"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.001275"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "-2.33%"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.004082"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.0001352"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "0.12%"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "0.86%"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.02573"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "-7.27%"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.05804"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "4.15%"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.01075"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "70.58%"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.007892"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "2.68%"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1417"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "-2.20%"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.008394"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "9.16%"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.007672"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "-5.01%"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.3366"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "-4.20%"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.00006995"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "2.28%"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "0.13%"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "-7.54%"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "0.10%"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.00002103"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "0.13%"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "0.13%"
